$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interest count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 955
$ws1.Range("F21").Value = 1180
$ws1.Range("F22").Value = 2859
$ws1.Range("F23").Value = 1405
$ws1.Range("F24").Value = 699
$ws1.Range("F31").Value = 593

# Sheet "全部类型" (All types) - update "想去人数" (interest count) column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 955
$ws4.Range("F33").Value = 1180
$ws4.Range("F34").Value = 2859
$ws4.Range("F35").Value = 1405
$ws4.Range("F36").Value = 699
$ws4.Range("F45").Value = 593
